# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15, 16) switch their table style from the
#    custom "Table_0" style to the built-in "Medium Style 2 - Accent 4"
#    style ({532C04C2-86F0-48D5-8A5A-A244AD2CA81F}).
# 2) The deck's active theme (the theme used by the slide master / all
#    slides) switches its 12 theme colors from the "Integral / Red Violet"
#    palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$newTableStyle = "{532C04C2-86F0-48D5-8A5A-A244AD2CA81F}"

foreach ($slideIdx in 14, 15, 16) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Swap the active theme's colour palette back to "Office" -----------
function Set-ThemeColor {
    param($ColorScheme, [int]$Index, [string]$Hex)

    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)

    # VBA-style RGB() packing: 0x00BBGGRR
    $ColorScheme.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

$tcs = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeColor $tcs 1  "000000"   # dk1
Set-ThemeColor $tcs 2  "FFFFFF"   # lt1
Set-ThemeColor $tcs 3  "44546A"   # dk2
Set-ThemeColor $tcs 4  "E7E6E6"   # lt2
Set-ThemeColor $tcs 5  "5B9BD5"   # accent1
Set-ThemeColor $tcs 6  "ED7D31"   # accent2
Set-ThemeColor $tcs 7  "A5A5A5"   # accent3
Set-ThemeColor $tcs 8  "FFC000"   # accent4
Set-ThemeColor $tcs 9  "4472C4"   # accent5
Set-ThemeColor $tcs 10 "70AD47"   # accent6
Set-ThemeColor $tcs 11 "0563C1"   # hlink
Set-ThemeColor $tcs 12 "954F72"   # folHlink
